# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (AC1, bold/
# centered/bordered style) onto the three new header cells so they match the
# look of the other headers in row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record is identical for every player row (2-46): 97 wins, 65 losses,
# 0 ties.
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 97   # column AD
    $ws.Cells.Item($row, 31).Value = 65   # column AE
    $ws.Cells.Item($row, 32).Value = 0    # column AF
}
